$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Insert a new column before column A to make room for the identifier column
$ws.Columns.Item(1).Insert()

# Header cell for the new column (bold, like the other header cells)
$ws.Range("A1").Value = "Identificador"
$ws.Range("A1").Font.Bold = $true

# Data cell for the new column
$ws.Range("A2").Value = 1

# Match the bestFit width Excel would compute for the new "Identificador" column
$ws.Columns.Item(1).ColumnWidth = 11.1666666666667
